$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the values in row 5 (B5:AH5) to 2 decimal places ("custom accuracy")
$row5 = @{
    "B5"  = 1.44
    "C5"  = 0.74
    "D5"  = 0.6
    "E5"  = 3.09
    "F5"  = 2.15
    "G5"  = 1.14
    "H5"  = 11.07
    "I5"  = 1.75
    "J5"  = 0.77
    "K5"  = 0.76
    "L5"  = 1.24
    "M5"  = 1.23
    "N5"  = 0.4
    "O5"  = 1.13
    "P5"  = 1.67
    "Q5"  = 1.21
    "R5"  = 0.68
    "S5"  = 0.28
    "T5"  = 9.99
    "U5"  = 3.68
    "V5"  = 1.04
    "W5"  = 2.41
    "X5"  = 1.09
    "Y5"  = 0.47
    "Z5"  = 4.9
    "AA5" = 0.92
    "AB5" = 0.97
    "AC5" = 1.11
    "AD5" = 1.24
    "AE5" = 0.56
    "AF5" = 10.65
    "AG5" = 0.46
    "AH5" = 1.31
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Remove the last data row (row 6) entirely, shrinking the used range to A1:AH5
$ws.Rows.Item(6).Delete()
